$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.030.40'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '3.320.51'
$ws.Range("E3").Value = '  +6.13%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.39%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.320.46'
$ws.Range("E8").Value = '  +6.46%  '
$ws.Range("E9").Value = '  +1.50%  '
$ws.Range("E10").Value = '  +3.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.75%  '
$ws.Range("E12").Value = '  +4.25%  '
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.47%  '
$ws.Range("D15").Value = '3.868.50'
$ws.Range("E15").Value = '  +6.14%  '
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '3.317.59'
$ws.Range("E17").Value = '  +6.05%  '
$ws.Range("D18").Value = '64.123.87'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '483.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("E22").Value = '  +6.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  +2.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.08%  '
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.33%  '
$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '29.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.107'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("E35").Value = '  +2.29%  '
$ws.Range("E36").Value = '  +3.72%  '
$ws.Range("D37").Value = '0.0₃0769'
$ws.Range("E37").Value = '  +8.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("E39").Value = '  +4.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '436.14'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("D41").Value = '3.065.86'
$ws.Range("E41").Value = '  +5.86%  '
$ws.Range("E42").Value = '  +3.83%  '
$ws.Range("E43").Value = '  +2.83%  '
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("E45").Value = '  +2.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.17%  '
$ws.Range("E48").Value = '  +15.07%  '
$ws.Range("E50").Value = '  +2.97%  '
$ws.Range("E51").Value = '  +2.01%  '